# Add a new "Financial Labels" sheet at the end of the workbook containing
# the label strings used for the financial-summary totals on the PMO
# Dashboard (Total Budget Approved / Total Actual Cost / Total Cost at
# Completion / Total Savings in PHP).

$wb = $excel.ActiveWorkbook

# --- create the new worksheet as the LAST tab ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Financial Labels"

# --- content --------------------------------------------------------------
$newSheet.Range("A1").Value = "Financial Labels"
$newSheet.Range("A2").Value = "Total Budget Approved"
$newSheet.Range("A3").Value = "Total Actual Cost"
$newSheet.Range("A4").Value = "Total Cost at Completion"
$newSheet.Range("A5").Value = "Total Savings in PHP"

# --- formatting: reuse the look of the other "label list" sheets ----------
# (bold header cell in row 1, plain body cells below it -- same visual
# style already used on the "Project Sizes" sheet).
$refSheet = $wb.Worksheets.Item("Project Sizes")

$refSheet.Range("A1").Copy() | Out-Null
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null

$refSheet.Range("A2").Copy() | Out-Null
$newSheet.Range("A2:A5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Make this new sheet the active / visible tab, as it is the sheet that was
# being worked on when the edit was made.
$newSheet.Activate()
